$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 198.26666
$ws.Range("I28").Value = 197.83333
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 197.83333
$ws.Range("L28").Value = 200
$ws.Range("M28").Value = 287.16667
$ws.Range("N28").Value = -1170
# Row 62
$ws.Range("H62").Value = 4167.8335
$ws.Range("I62").Value = 5901.6665
$ws.Range("J62").Value = 2434
$ws.Range("K62").Value = 5901.6665
$ws.Range("L62").Value = 2434
$ws.Range("M62").Value = -5277.6665
$ws.Range("N62").Value = -3682
# Row 65
$ws.Range("H65").Value = 4167.8335
$ws.Range("I65").Value = 5901.6665
$ws.Range("J65").Value = 2434
$ws.Range("K65").Value = 29508.3325
$ws.Range("L65").Value = 12170
$ws.Range("M65").Value = -26388.3325
$ws.Range("N65").Value = -18410
# Row 98
$ws.Range("H98").Value = 2445.5557
$ws.Range("I98").Value = 2445.5557
$ws.Range("K98").Value = 2445.5557
$ws.Range("M98").Value = -947.5556999999999
# Row 107
$ws.Range("H107").Value = 17858030
$ws.Range("I107").Value = 19231340
$ws.Range("K107").Value = 19231340
$ws.Range("M107").Value = -19229420
# Row 113
$ws.Range("H113").Value = 2477.3125
$ws.Range("I113").Value = 2348.6667
$ws.Range("J113").Value = 2642.7144
$ws.Range("K113").Value = 2348.6667
$ws.Range("L113").Value = 2642.7144
$ws.Range("M113").Value = 905.3332999999998
$ws.Range("N113").Value = -9150.714400000001
# Row 122
$ws.Range("H122").Value = 2445.5557
$ws.Range("I122").Value = 2445.5557
$ws.Range("K122").Value = 7336.6671
$ws.Range("M122").Value = -4886.6671
# Row 129
$ws.Range("H129").Value = 1179.4286
$ws.Range("I129").Value = 703.75
$ws.Range("J129").Value = 1369.7
$ws.Range("K129").Value = 2111.25
$ws.Range("L129").Value = 4109.1
$ws.Range("M129").Value = 2888.75
$ws.Range("N129").Value = -14109.1

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12212.186
$ws.Range("I32").Value = 9608.431
$ws.Range("J32").Value = 18778.174
$ws.Range("K32").Value = 9608.431
$ws.Range("L32").Value = 18778.174
$ws.Range("M32").Value = -9321.431
$ws.Range("N32").Value = -19352.174
# Row 61
$ws.Range("H61").Value = 300167.38
$ws.Range("I61").Value = 7683.45
$ws.Range("K61").Value = 7683.45
$ws.Range("M61").Value = -7471.45
# Row 74
$ws.Range("H74").Value = 1760.697
$ws.Range("I74").Value = 1641.1765
$ws.Range("J74").Value = 1887.6875
$ws.Range("K74").Value = 1641.1765
$ws.Range("L74").Value = 1887.6875
$ws.Range("M74").Value = -767.1765
$ws.Range("N74").Value = -3635.6875
# Row 77
$ws.Range("H77").Value = 1760.697
$ws.Range("I77").Value = 1641.1765
$ws.Range("J77").Value = 1887.6875
$ws.Range("K77").Value = 8205.8825
$ws.Range("L77").Value = 9438.4375
$ws.Range("M77").Value = -3837.8825
$ws.Range("N77").Value = -18174.4375
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 125
$ws.Range("H125").Value = 77111
$ws.Range("J125").Value = 77111
$ws.Range("L125").Value = 77111
$ws.Range("N125").Value = -86951
# Row 136
$ws.Range("H136").Value = 300167.38
$ws.Range("I136").Value = 7683.45
$ws.Range("K136").Value = 23050.35
$ws.Range("M136").Value = -20500.35

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 202084.67
$ws.Range("I107").Value = 202084.67
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 202084.67
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -200164.67
$ws.Range("N107").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 53334.332
$ws.Range("I4").Value = 3000
$ws.Range("J4").Value = 63401.2
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 63401.2
$ws.Range("M4").Value = -2888
$ws.Range("N4").Value = -63625.2
# Row 16
$ws.Range("H16").Value = 1877.3572
$ws.Range("I16").Value = 1603.5
$ws.Range("J16").Value = 2082.75
$ws.Range("K16").Value = 1603.5
$ws.Range("L16").Value = 2082.75
$ws.Range("M16").Value = -1316.5
$ws.Range("N16").Value = -2656.75
# Row 99
$ws.Range("H99").Value = 5427.304
$ws.Range("I99").Value = 6772.2354
$ws.Range("J99").Value = 1616.6666
$ws.Range("K99").Value = 6772.2354
$ws.Range("L99").Value = 1616.6666
$ws.Range("M99").Value = -5274.2354
$ws.Range("N99").Value = -4612.6666
# Row 113
$ws.Range("H113").Value = 1877.3572
$ws.Range("I113").Value = 1603.5
$ws.Range("J113").Value = 2082.75
$ws.Range("K113").Value = 1603.5
$ws.Range("L113").Value = 2082.75
$ws.Range("M113").Value = 566.5
$ws.Range("N113").Value = -6422.75
# Row 126
$ws.Range("H126").Value = 5427.304
$ws.Range("I126").Value = 6772.2354
$ws.Range("J126").Value = 1616.6666
$ws.Range("K126").Value = 20316.7062
$ws.Range("L126").Value = 4849.9998
$ws.Range("M126").Value = -17846.7062
$ws.Range("N126").Value = -9789.9998

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 643100
$ws.Range("I4").Value = 643100
$ws.Range("K4").Value = 1929300
$ws.Range("M4").Value = -1929188
# Row 22
$ws.Range("H22").Value = 2142.5715
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2142.5715
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6427.7145
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6765.7145
# Row 27
$ws.Range("H27").Value = 2142.5715
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2142.5715
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6427.7145
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -6631.7145
# Row 32
$ws.Range("H32").Value = 5333.3335
$ws.Range("J32").Value = 5333.3335
$ws.Range("L32").Value = 16000.0005
$ws.Range("N32").Value = -16566.0005
# Row 86
$ws.Range("H86").Value = 850
$ws.Range("I86").Value = 850
$ws.Range("K86").Value = 2550
$ws.Range("M86").Value = -1364
# Row 89
$ws.Range("H89").Value = 850
$ws.Range("I89").Value = 850
$ws.Range("K89").Value = 7650
$ws.Range("M89").Value = -1722
# Row 122
$ws.Range("H122").Value = 2964.9788
$ws.Range("I122").Value = 661.9231
$ws.Range("J122").Value = 3845.5588
$ws.Range("K122").Value = 5957.3079
$ws.Range("L122").Value = 34610.0292
$ws.Range("M122").Value = -3507.3079
$ws.Range("N122").Value = -39510.0292
# Row 132
$ws.Range("H132").Value = 2802821
$ws.Range("I132").Value = 852
$ws.Range("J132").Value = 3923608.5
$ws.Range("K132").Value = 7668
$ws.Range("L132").Value = 35312476.5
$ws.Range("M132").Value = -5138
$ws.Range("N132").Value = -35317536.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 80
$ws.Range("H80").Value = 7125.25
$ws.Range("I80").Value = 7982.647
$ws.Range("J80").Value = 2266.6667
$ws.Range("K80").Value = 7982.647
$ws.Range("L80").Value = 2266.6667
$ws.Range("M80").Value = -6984.647
$ws.Range("N80").Value = -4262.6667
# Row 83
$ws.Range("H83").Value = 7125.25
$ws.Range("I83").Value = 7982.647
$ws.Range("J83").Value = 2266.6667
$ws.Range("K83").Value = 39913.235
$ws.Range("L83").Value = 11333.3335
$ws.Range("M83").Value = -34921.235
$ws.Range("N83").Value = -21317.3335
# Row 102
$ws.Range("H102").Value = 2059.3635
$ws.Range("I102").Value = 1821.0646
$ws.Range("J102").Value = 2627.6155
$ws.Range("K102").Value = 1821.0646
$ws.Range("L102").Value = 2627.6155
$ws.Range("M102").Value = -199.0645999999999
$ws.Range("N102").Value = -5871.6155
# Row 122
$ws.Range("H122").Value = 8288041.5
$ws.Range("I122").Value = 3243249.8
$ws.Range("J122").Value = 13893366
$ws.Range("K122").Value = 9729749.399999999
$ws.Range("L122").Value = 41680098
$ws.Range("M122").Value = -9727299.399999999
$ws.Range("N122").Value = -41684998
# Row 126
$ws.Range("H126").Value = 13831.25
$ws.Range("I126").Value = 16215.385
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 48646.155
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -46176.155
$ws.Range("N126").Value = -15440
# Row 132
$ws.Range("H132").Value = 4527.543
$ws.Range("I132").Value = 4224.6523
$ws.Range("J132").Value = 5108.0835
$ws.Range("K132").Value = 12673.9569
$ws.Range("L132").Value = 15324.2505
$ws.Range("M132").Value = -10143.9569
$ws.Range("N132").Value = -20384.2505
# Row 135
$ws.Range("H135").Value = 56780
$ws.Range("J135").Value = 56780
$ws.Range("L135").Value = 56780
$ws.Range("N135").Value = -66920
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 140
$ws.Range("H140").Value = 31025.191
$ws.Range("J140").Value = 31025.191
$ws.Range("L140").Value = 31025.191
$ws.Range("N140").Value = -41385.191
# Row 141
$ws.Range("H141").Value = 69349
$ws.Range("J141").Value = 69349
$ws.Range("L141").Value = 69349
$ws.Range("N141").Value = -79709

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2095.8333
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 2180.9524
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2180.9524
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -2770.9524
# Row 27
$ws.Range("H27").Value = 2095.8333
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 2180.9524
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 2180.9524
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -2394.9524
# Row 132
$ws.Range("H132").Value = 9265145
$ws.Range("I132").Value = 11500846
$ws.Range("J132").Value = 2957
$ws.Range("K132").Value = 34502538
$ws.Range("L132").Value = 8871
$ws.Range("M132").Value = -34500008
$ws.Range("N132").Value = -13931

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 100003
$ws.Range("J2").Value = 100003
$ws.Range("L2").Value = 100003
$ws.Range("N2").Value = -100227
# Row 100
$ws.Range("H100").Value = 416.66666
$ws.Range("I100").Value = 416.66666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 833.33332
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -292.33332
$ws.Range("N100").ClearContents()
# Row 126
$ws.Range("H126").Value = 1213.8
$ws.Range("I126").Value = 927.9091
$ws.Range("K126").Value = 2783.7273
$ws.Range("M126").Value = -313.7273
# Row 132
$ws.Range("H132").Value = 1220.5151
$ws.Range("I132").Value = 886.16394
$ws.Range("J132").Value = 5299.6
$ws.Range("K132").Value = 2658.49182
$ws.Range("L132").Value = 15898.8
$ws.Range("M132").Value = -128.4918200000002
$ws.Range("N132").Value = -20958.8
